# courrier_radiation.template.docx — fix(docs): correction des codes postaux
#
# Three paragraphs hold "{STRUCTURE_VILLE}, {STRUCTURE_CODE_POSTAL}" (two short
# "address block" paragraphs, one long sentence). The city/postal-code order was
# wrong; swap to "{STRUCTURE_CODE_POSTAL}, {STRUCTURE_VILLE}" everywhere.

$d = $word.ActiveDocument

$oldShort = "{STRUCTURE_VILLE}, {STRUCTURE_CODE_POSTAL}"
$newShort = "{STRUCTURE_CODE_POSTAL}, {STRUCTURE_VILLE}"

$oldLong = "{STRUCTURE_ADRESSE} {STRUCTURE_VILLE}, {STRUCTURE_CODE_POSTAL} pour"
$newLong = "{STRUCTURE_ADRESSE} {STRUCTURE_CODE_POSTAL}, {STRUCTURE_VILLE} pour"

$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $para = $d.Paragraphs($i)
    $r = $para.Range

    if ($r.Text.Contains($oldShort)) {
        $r.Find.Execute($oldShort, $false, $false, $false, $false, $false, $true, 0, $false, $newShort, 2)
    }
    elseif ($r.Text.Contains($oldLong)) {
        $r.Find.Execute($oldLong, $false, $false, $false, $false, $false, $true, 0, $false, $newLong, 2)
    }
}
